$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.ClearFormats()
}

Set-TextValue $ws.Range("D2") "26.486.58"
Set-TextValue $ws.Range("E2") "  -3.25%  "

Set-TextValue $ws.Range("D3") "1.804.14"
Set-TextValue $ws.Range("E3") "  -3.07%  "

Set-TextValue $ws.Range("D4") "1.006"
Set-TextValue $ws.Range("E4") "  +0.27%  "

Set-TextValue $ws.Range("D5") "1.006"
Set-TextValue $ws.Range("E5") "  +0.36%  "

Set-TextValue $ws.Range("D6") "308.04"
Set-TextValue $ws.Range("E6") "  -2.33%  "

Set-TextValue $ws.Range("D7") "0.4529"
Set-TextValue $ws.Range("E7") "  -1.92%  "

Set-TextValue $ws.Range("D8") "0.3655"
Set-TextValue $ws.Range("E8") "  -1.64%  "

Set-TextValue $ws.Range("D9") "0.07099"
Set-TextValue $ws.Range("E9") "  -2.99%  "

Set-TextValue $ws.Range("D10") "0.8730"
Set-TextValue $ws.Range("E10") "  -1.80%  "

Set-TextValue $ws.Range("D11") "0.07785"
Set-TextValue $ws.Range("E11") "  -0.63%  "

Set-TextValue $ws.Range("D12") "19.27"
Set-TextValue $ws.Range("E12") "  -3.80%  "

Set-TextValue $ws.Range("D13") "1.852.56"
Set-TextValue $ws.Range("E13") "  -3.50%  "

Set-TextValue $ws.Range("D14") "5.269"
Set-TextValue $ws.Range("E14") "  -2.34%  "

Set-TextValue $ws.Range("D15") "6.332"

Set-TextValue $ws.Range("D16") "86.44"
Set-TextValue $ws.Range("E16") "  -5.81%  "

Set-TextValue $ws.Range("E17") "  +0.39%  "

Set-TextValue $ws.Range("D18") "0.000008543"
Set-TextValue $ws.Range("E18") "  -4.65%  "

Set-TextValue $ws.Range("D19") "1.006"
Set-TextValue $ws.Range("E19") "  +0.40%  "

Set-TextValue $ws.Range("D20") "26.529.17"
Set-TextValue $ws.Range("E20") "  -3.12%  "

Set-TextValue $ws.Range("D21") "14.19"
Set-TextValue $ws.Range("E21") "  -4.12%  "

Set-TextValue $ws.Range("E22") "  -3.29%  "

Set-TextValue $ws.Range("D23") "2.061.01"
Set-TextValue $ws.Range("E23") "  +0.46%  "

Set-TextValue $ws.Range("E24") "  -1.88%  "

Set-TextValue $ws.Range("D25") "1.979"
Set-TextValue $ws.Range("E25") "  +2.06%  "

Set-TextValue $ws.Range("D26") "150.42"
Set-TextValue $ws.Range("E26") "  -1.07%  "

Set-TextValue $ws.Range("D27") "17.84"
Set-TextValue $ws.Range("E27") "  -3.11%  "

Set-TextValue $ws.Range("D28") "1.995"
Set-TextValue $ws.Range("E28") "  -2.81%  "

Set-TextValue $ws.Range("D29") "113.21"
Set-TextValue $ws.Range("E29") "  -2.62%  "

Set-TextValue $ws.Range("D30") "4.876"
Set-TextValue $ws.Range("E30") "  -4.44%  "

Set-TextValue $ws.Range("D31") "0.08653"
Set-TextValue $ws.Range("E31") "  -2.16%  "

Set-TextValue $ws.Range("D32") "3.113"
Set-TextValue $ws.Range("E32") "  +0.08%  "

Set-TextValue $ws.Range("D33") "0.7284"
Set-TextValue $ws.Range("E33") "  -4.83%  "

Set-TextValue $ws.Range("E34") "  -1.74%  "

Set-TextValue $ws.Range("E35") "  -5.21%  "

Set-TextValue $ws.Range("B36") "RenderToken"
Set-TextValue $ws.Range("C36") "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue $ws.Range("D36") "2.525"
Set-TextValue $ws.Range("E36") "  -6.77%  "

Set-TextValue $ws.Range("B37") "TrustWalletToken"
Set-TextValue $ws.Range("C37") "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-TextValue $ws.Range("D37") "1.078"
Set-TextValue $ws.Range("E37") "  -0.22%  "

Set-TextValue $ws.Range("B38") "VeChain"
Set-TextValue $ws.Range("C38") "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue $ws.Range("D38") "0.01910"
Set-TextValue $ws.Range("E38") "  -2.38%  "

Set-TextValue $ws.Range("B39") "Hedera"
Set-TextValue $ws.Range("C39") "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue $ws.Range("D39") "0.05074"
Set-TextValue $ws.Range("E39") "  -3.23%  "

Set-TextValue $ws.Range("B40") "MXToken"
Set-TextValue $ws.Range("C40") "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextValue $ws.Range("D40") "2.863"
Set-TextValue $ws.Range("E40") "  -4.24%  "

Set-TextValue $ws.Range("B41") "FraxShare"
Set-TextValue $ws.Range("C41") "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextValue $ws.Range("D41") "6.899"
Set-TextValue $ws.Range("E41") "  -2.41%  "

Set-TextValue $ws.Range("B42") "TheSandbox"
Set-TextValue $ws.Range("C42") "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
Set-TextValue $ws.Range("D42") "0.4906"
Set-TextValue $ws.Range("E42") "  -4.67%  "

Set-TextValue $ws.Range("B43") "Algorand"
Set-TextValue $ws.Range("C43") "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-TextValue $ws.Range("D43") "0.1568"
Set-TextValue $ws.Range("E43") "  -4.81%  "

Set-TextValue $ws.Range("B44") "Aptos"
Set-TextValue $ws.Range("C44") "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextValue $ws.Range("D44") "8.130"
Set-TextValue $ws.Range("E44") "  -3.37%  "

Set-TextValue $ws.Range("B45") "PaxDollar"
Set-TextValue $ws.Range("C45") "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
Set-TextValue $ws.Range("D45") "1.006"
Set-TextValue $ws.Range("E45") "  +0.45%  "

Set-TextValue $ws.Range("B46") "Decentraland"
Set-TextValue $ws.Range("C46") "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
Set-TextValue $ws.Range("D46") "0.4595"
Set-TextValue $ws.Range("E46") "  -4.48%  "

Set-TextValue $ws.Range("B47") "Quant"
Set-TextValue $ws.Range("C47") "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
Set-TextValue $ws.Range("D47") "101.50"
Set-TextValue $ws.Range("E47") "  -1.25%  "

Set-TextValue $ws.Range("B48") "EnergySwap"
Set-TextValue $ws.Range("C48") "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue $ws.Range("D48") "9.916"
Set-TextValue $ws.Range("E48") "  -4.36%  "

Set-TextValue $ws.Range("B49") "NEARProtocol"
Set-TextValue $ws.Range("C49") "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue $ws.Range("D49") "1.578"
Set-TextValue $ws.Range("E49") "  -4.18%  "

Set-TextValue $ws.Range("B50") "Cronos"
Set-TextValue $ws.Range("C50") "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextValue $ws.Range("D50") "0.05993"
Set-TextValue $ws.Range("E50") "  -3.67%  "

Set-TextValue $ws.Range("B51") "Aave"
Set-TextValue $ws.Range("C51") "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextValue $ws.Range("D51") "63.35"
Set-TextValue $ws.Range("E51") "  -3.19%  "
